$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3825411565577213
$ws.Range("C2").Value = 0.06080760717030387
$ws.Range("D2").Value = 0.1972768006716592
$ws.Range("E2").Value = 0.1687696096472209
$ws.Range("F2").Value = 1.424585000670405
$ws.Range("I2").Value = 0.7082064001760067
$ws.Range("J2").Value = 0.1873394869776064
$ws.Range("K2").Value = 0.416912196189628
$ws.Range("O2").Value = 3.444437899617611
$ws.Range("B3").Value = 0.3426740460886037
$ws.Range("C3").Value = 0.05361555193017864
$ws.Range("D3").Value = 0.1905918361528904
$ws.Range("E3").Value = 0.1643571820370866
$ws.Range("F3").Value = 1.429303185283011
$ws.Range("I3").Value = 0.7158796129890419
$ws.Range("J3").Value = 0.1834558576620609
$ws.Range("K3").Value = 0.3719707251526927
$ws.Range("O3").Value = 3.469712867894444
$ws.Range("B4").Value = 0.3182053267265985
$ws.Range("C4").Value = 0.04919238347139299
$ws.Range("D4").Value = 0.1865624613252663
$ws.Range("E4").Value = 0.1617301785162901
$ws.Range("F4").Value = 1.432980918564809
$ws.Range("I4").Value = 0.7209753403111172
$ws.Range("J4").Value = 0.1811772573572625
$ws.Range("K4").Value = 0.3443737628843166
$ws.Range("O4").Value = 3.487160028937922
$ws.Range("B5").Value = 0.308237226976189
$ws.Range("C5").Value = 0.04738818196059924
$ws.Range("D5").Value = 0.1849394946411138
$ws.Range("E5").Value = 0.1606803902649645
$ws.Range("F5").Value = 1.434675993104449
$ws.Range("I5").Value = 0.7231485015044541
$ws.Range("J5").Value = 0.1802753754285433
$ws.Range("K5").Value = 0.3331277659757745
$ws.Range("O5").Value = 3.494754645539984
$ws.Range("B6").Value = 0.3065822374375671
$ws.Range("C6").Value = 0.04708849419462524
$ws.Range("D6").Value = 0.1846711555831888
$ws.Range("E6").Value = 0.1605073278374469
$ws.Range("F6").Value = 1.434969320372645
$ws.Range("I6").Value = 0.7235151872924028
$ws.Range("J6").Value = 0.1801272298087042
$ws.Range("K6").Value = 0.3312603959683997
$ws.Range("O6").Value = 3.496045000701216
$ws.Range("B7").Value = 0.3180708800554441
$ws.Range("C7").Value = 0.04916805823923198
$ws.Range("D7").Value = 0.1865404961896218
$ws.Range("E7").Value = 0.1617159366669298
$ws.Range("F7").Value = 1.433002983789535
$ws.Range("I7").Value = 0.7210042571966824
$ws.Range("J7").Value = 0.1811649862640579
$ws.Range("K7").Value = 0.3442220945863994
$ws.Range("O7").Value = 3.48726048998428
$ws.Range("B8").Value = 0.3687932953046982
$ws.Range("C8").Value = 0.0583293364677786
$ws.Range("D8").Value = 0.1949562833189304
$ws.Range("E8").Value = 0.1672311664672392
$ws.Range("F8").Value = 1.426049809422501
$ws.Range("I8").Value = 0.7107723325058544
$ws.Range("J8").Value = 0.1859784295885873
$ws.Range("K8").Value = 0.4014173298091634
$ws.Range("O8").Value = 3.452752435521049
$ws.Range("B9").Value = 0.4683154017182574
$ws.Range("C9").Value = 0.07623421192150204
$ws.Range("D9").Value = 0.2120521081270965
$ws.Range("E9").Value = 0.1786973675917238
$ws.Range("F9").Value = 1.418608972497921
$ws.Range("I9").Value = 0.6937587135356829
$ws.Range("J9").Value = 0.1962581699804531
$ws.Range("K9").Value = 0.5135314631701249
$ws.Range("O9").Value = 3.400388407791922
$ws.Range("B10").Value = 0.5414452157916969
$ws.Range("C10").Value = 0.08934925765666435
$ws.Range("D10").Value = 0.2249692238504508
$ws.Range("E10").Value = 0.1875171254761767
$ws.Range("F10").Value = 1.416919572574898
$ws.Range("I10").Value = 0.6831209828711877
$ws.Range("J10").Value = 0.2043239766557434
$ws.Range("K10").Value = 0.5958499032641669
$ws.Range("O10").Value = 3.371258813343502
$ws.Range("B11").Value = 0.5747118371222655
$ws.Range("C11").Value = 0.09530649450982764
$ws.Range("D11").Value = 0.2309221928581593
$ws.Range("E11").Value = 0.1916151523982634
$ws.Range("F11").Value = 1.416971627563683
$ws.Range("I11").Value = 0.6786864215858657
$ws.Range("J11").Value = 0.2081050173505901
$ws.Range("K11").Value = 0.6332828751882289
$ws.Range("O11").Value = 3.360038181655483
$ws.Range("B12").Value = 0.5873084391926682
$ws.Range("C12").Value = 0.09756099764879878
$ws.Range("D12").Value = 0.2331873837916874
$ws.Range("E12").Value = 0.193179276325445
$ws.Range("F12").Value = 1.41710933514095
$ws.Range("I12").Value = 0.6770653943257336
$ws.Range("J12").Value = 0.2095528808557674
$ws.Range("K12").Value = 0.6474551752294531
$ws.Range("O12").Value = 3.356081404844645
$ws.Range("B13").Value = 0.5845955794102338
$ws.Range("C13").Value = 0.09707551269804071
$ws.Range("D13").Value = 0.2326990502921547
$ws.Range("E13").Value = 0.1928418685836135
$ws.Range("F13").Value = 1.417074429505277
$ws.Range("I13").Value = 0.6774119200686606
$ws.Range("J13").Value = 0.2092403433599941
$ws.Range("K13").Value = 0.644403050865094
$ws.Range("O13").Value = 3.356920566584449
$ws.Range("B14").Value = 0.575748185945713
$ws.Range("C14").Value = 0.09549200192358853
$ws.Range("D14").Value = 0.2311083331460253
$ws.Range("E14").Value = 0.1917435878402429
$ws.Range("F14").Value = 1.416980592056476
$ws.Range("I14").Value = 0.6785518910149761
$ws.Range("J14").Value = 0.2082238120957953
$ws.Range("K14").Value = 0.6344488981885661
$ws.Range("O14").Value = 3.359706797094418
$ws.Range("B15").Value = 0.5703287849520677
$ws.Range("C15").Value = 0.09452187410255419
$ws.Range("D15").Value = 0.2301353926581697
$ws.Range("E15").Value = 0.1910724579593861
$ws.Range("F15").Value = 1.416938480339724
$ws.Range("I15").Value = 0.6792577432472733
$ws.Range("J15").Value = 0.2076032489247126
$ws.Range("K15").Value = 0.6283513120833106
$ws.Range("O15").Value = 3.361451508176998
$ws.Range("B16").Value = 0.5392710961704381
$ws.Range("C16").Value = 0.08895975142371526
$ws.Range("D16").Value = 0.2245817199272437
$ws.Range("E16").Value = 0.1872510327935544
$ws.Range("F16").Value = 1.416932682528852
$ws.Range("I16").Value = 0.6834189388949525
$ws.Range("J16").Value = 0.2040791260096739
$ws.Range("K16").Value = 0.593403224627167
$ws.Range("O16").Value = 3.372032988741381
$ws.Range("B17").Value = 0.5202176022428944
$ws.Range("C17").Value = 0.08554522890648286
$ws.Range("D17").Value = 0.2211943281627811
$ws.Range("E17").Value = 0.1849286651853177
$ws.Range("F17").Value = 1.417139292860298
$ws.Range("I17").Value = 0.6860753702278473
$ws.Range("J17").Value = 0.2019458274802872
$ws.Range("K17").Value = 0.5719595759017864
$ws.Range("O17").Value = 3.379044632294892
$ws.Range("B18").Value = 0.5092585297303174
$ws.Range("C18").Value = 0.08358045748417453
$ws.Range("D18").Value = 0.2192532366667592
$ws.Range("E18").Value = 0.1836009872254252
$ws.Range("F18").Value = 1.417335365879367
$ws.Range("I18").Value = 0.6876413578845089
$ws.Range("J18").Value = 0.2007293410867703
$ws.Range("K18").Value = 0.5596244724696646
$ws.Range("O18").Value = 3.383268656509699
$ws.Range("B19").Value = 0.5055479954968405
$ws.Range("C19").Value = 0.08291508037842732
$ws.Range("D19").Value = 0.2185972649758128
$ws.Range("E19").Value = 0.1831528489213383
$ws.Range("F19").Value = 1.417415018271171
$ws.Range("I19").Value = 0.6881781124719133
$ws.Range("J19").Value = 0.2003192689392108
$ws.Range("K19").Value = 0.55544782007712
$ws.Range("O19").Value = 3.384731654931556
$ws.Range("B20").Value = 0.5222458847871394
$ws.Range("C20").Value = 0.08590879710740751
$ws.Range("D20").Value = 0.2215541727290997
$ws.Range("E20").Value = 0.1851750486876682
$ws.Range("F20").Value = 1.417109305186983
$ws.Range("I20").Value = 0.6857886471155297
$ws.Range("J20").Value = 0.202171831032274
$ws.Range("K20").Value = 0.5742424270614777
$ws.Range("O20").Value = 3.378278449716646
$ws.Range("B21").Value = 0.5783469054586021
$ws.Range("C21").Value = 0.09595715547331451
$ws.Range("D21").Value = 0.2315752696982969
$ws.Range("E21").Value = 0.1920658464365346
$ws.Range("F21").Value = 1.417004952043968
$ws.Range("I21").Value = 0.6782154726525285
$ws.Range("J21").Value = 0.2085219561230218
$ws.Range("K21").Value = 0.6373727546415751
$ws.Range("O21").Value = 3.358880480632251
$ws.Range("B22").Value = 0.6150075096577723
$ws.Range("C22").Value = 0.1025162891919251
$ws.Range("D22").Value = 0.2381882891205009
$ws.Range("E22").Value = 0.1966409896936625
$ws.Range("F22").Value = 1.417624517833787
$ws.Range("I22").Value = 0.6736054908103561
$ws.Range("J22").Value = 0.2127657619222987
$ws.Range("K22").Value = 0.6786156854005299
$ws.Range("O22").Value = 3.347906253817541
$ws.Range("B23").Value = 0.5954417032301649
$ws.Range("C23").Value = 0.09901632508413627
$ws.Range("D23").Value = 0.234653014381621
$ws.Range("E23").Value = 0.1941926168658128
$ws.Range("F23").Value = 1.417230913998893
$ws.Range("I23").Value = 0.6760348369953775
$ws.Range("J23").Value = 0.2104922024532385
$ws.Range("K23").Value = 0.656605295878478
$ws.Range("O23").Value = 3.353607458916315
$ws.Range("B24").Value = 0.5213289136647745
$ws.Range("C24").Value = 0.0857444332653472
$ws.Range("D24").Value = 0.2213914671650627
$ws.Range("E24").Value = 0.1850636353966877
$ws.Range("F24").Value = 1.417122621859633
$ws.Range("I24").Value = 0.6859181538724393
$ws.Range("J24").Value = 0.202069623753431
$ws.Range("K24").Value = 0.5732103713796448
$ws.Range("O24").Value = 3.378624239886364
$ws.Range("B25").Value = 0.441388468098296
$ws.Range("C25").Value = 0.0713972291537317
$ws.Range("D25").Value = 0.2073642579644996
$ws.Range("E25").Value = 0.1755259054453191
$ws.Range("F25").Value = 1.419958649565984
$ws.Range("I25").Value = 0.6980344832244008
$ws.Range("J25").Value = 0.1933871493188377
$ws.Range("K25").Value = 0.4832089932229167
$ws.Range("O25").Value = 3.412914473483568
